$d = $word.ActiveDocument

# The document contains a short "Isaiah" paragraph (body text, not a
# heading) that sits directly after the "ISA" Heading2 paragraph. It
# should be removed entirely, including its paragraph mark, so the
# "ISA" heading paragraph becomes immediately followed by the paragraph
# that previously came after the "Isaiah" paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Isaiah" -and $p.Style.NameLocal -ne "Heading 2") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
